$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.541.43'
$ws.Range('E2').Value = '  -0.31%  '

$ws.Range('D3').Value = '1.626.61'
$ws.Range('E3').Value = '  -0.58%  '

$ws.Range('E4').Value = '  -0.17%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.13'
$ws.Range('E5').Value = '  -0.20%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.520'
$ws.Range('E6').Value = '  -0.41%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.19%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '23.39'
$ws.Range('E8').Value = '  +1.97%  '

$ws.Range('E9').Value = '  +2.21%  '

$ws.Range('E10').Value = '  +0.22%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0879'
$ws.Range('E11').Value = '  -1.73%  '

$ws.Range('D12').Value = '1.854.40'
$ws.Range('E12').Value = '  -0.77%  '

$ws.Range('D13').Value = '1.614.21'
$ws.Range('E13').Value = '  -3.18%  '

$ws.Range('E14').Value = '  +0.35%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.554'
$ws.Range('E15').Value = '  -1.06%  '

$ws.Range('E16').Value = '  +1.41%  '

$ws.Range('D17').Value = '27.505.56'
$ws.Range('E17').Value = '  -0.44%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '230.30'
$ws.Range('E18').Value = '  +0.39%  '

$ws.Range('E19').Value = '  -0.58%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.57'
$ws.Range('E20').Value = '  -2.10%  '

$ws.Range('E21').Value = '  -0.03%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.48'
$ws.Range('E22').Value = '  +4.59%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.37'
$ws.Range('E23').Value = '  +2.08%  '

$ws.Range('E24').Value = '  +8.55%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '149.12'
$ws.Range('E25').Value = '  -0.73%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.89'
$ws.Range('E26').Value = '  -0.50%  '

$ws.Range('E27').Value = '  +0.08%  '

$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.54'
$ws.Range('E28').Value = '  -0.49%  '

$ws.Range('B29').Value = 'BinanceUSD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  -0.20%  '

$ws.Range('E30').Value = '  -0.46%  '

$ws.Range('E31').Value = '  -0.19%  '

$ws.Range('E32').Value = '  -0.55%  '

$ws.Range('D33').Value = '1.468.83'
$ws.Range('E33').Value = '  +1.07%  '

$ws.Range('E34').Value = '  -1.51%  '

$ws.Range('E35').Value = '  -0.90%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.33'
$ws.Range('E36').Value = '  -1.90%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.943'
$ws.Range('E37').Value = '  +5.44%  '

$ws.Range('E38').Value = '  +0.67%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.876'
$ws.Range('E39').Value = '  +0.20%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.555'
$ws.Range('E40').Value = '  -1.67%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.04'
$ws.Range('E41').Value = '  +2.09%  '

$ws.Range('E42').Value = '  -0.14%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '67.94'
$ws.Range('E43').Value = '  -2.61%  '

$ws.Range('E44').Value = '  +0.59%  '

$ws.Range('E45').Value = '  -1.53%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '5.34'
$ws.Range('E46').Value = '  -4.62%  '

$ws.Range('E47').Value = '  +3.79%  '

$ws.Range('D48').Value = '1.764.39'
$ws.Range('E48').Value = '  -0.87%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '87.36'

$ws.Range('E50').Value = '  -0.82%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0995'
$ws.Range('E51').Value = '  +0.96%  '
